$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Give the new headers the same look as the existing header row (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for rows 2-49: column I values, then column J values
$iValues = @(4,6,7,7,6,7,7,1,6,7,5,5,5,8,8,4,4,5,7,11,7,7,8,1,8,4,5,5,6,7,8,4,6,8,7,7,9,5,9,5,7,5,9,5,5,9,6,8)
$jValues = @(4,6,7,7,6,7,8,1,7,7,5,5,5,8,9,5,5,5,7,11,7,7,8,1,8,5,6,5,7,7,8,5,6,8,8,7,9,6,9,5,7,6,9,5,5,9,6,8)

for ($row = 2; $row -le 49; $row++) {
    $idx = $row - 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
